$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("conversions")

# Header row: new columns J:M (Goals%, Assist%, PKatt%, Touch%) plus shifted headers N:Y
$ws.Range("J1").Value = "Goals%"
$ws.Range("K1").Value = "Assist%"
$ws.Range("L1").Value = "PKatt%"
$ws.Range("M1").Value = "Touch%"
$ws.Range("N1").Value = "PrgC"
$ws.Range("O1").Value = "Tkl"
$ws.Range("P1").Value = "TklW"
$ws.Range("Q1").Value = "blkSh"
$ws.Range("R1").Value = "blkPass"
$ws.Range("S1").Value = "Int"
$ws.Range("T1").Value = "Clr"
$ws.Range("U1").Value = "Err"
$ws.Range("V1").Value = "Fls"
$ws.Range("W1").Value = "Fld"
$ws.Range("X1").Value = "CrdY"
$ws.Range("Y1").Value = "CrdR"

# Recalibrated data values for rows 2-17, columns B through Y
# row 2
$ws.Range("B2").Value = [double]"0.164667197464969"
$ws.Range("C2").Value = [double]"0.217201918284426"
$ws.Range("D2").Value = [double]"0.17590169000081299"
$ws.Range("E2").Value = [double]"0.40520447199689302"
$ws.Range("F2").Value = [double]"0.20668412740338901"
$ws.Range("G2").Value = [double]"1.50451346496912E-2"
$ws.Range("H2").Value = [double]"0.34616988421265399"
$ws.Range("I2").Value = [double]"0.253445077426668"
$ws.Range("J2").Value = [double]"0.379275369550862"
$ws.Range("K2").Value = [double]"0.20003367860400101"
$ws.Range("L2").Value = [double]"0.79062889547406201"
$ws.Range("M2").Value = [double]"1.13531966159157E-2"
$ws.Range("N2").Value = [double]"0.33807672736913902"
$ws.Range("O2").Value = [double]"0.16600524746682199"
$ws.Range("P2").Value = [double]"0.29821886500973999"
$ws.Range("Q2").Value = [double]"0.163530412494822"
$ws.Range("R2").Value = [double]"7.7219274401898894E-2"
$ws.Range("S2").Value = [double]"0.18040985685611699"
$ws.Range("T2").Value = [double]"0.33161231288762499"
$ws.Range("U2").Value = [double]"0.64297012785460494"
$ws.Range("V2").Value = [double]"0.27143847761553203"
$ws.Range("W2").Value = [double]"0.16739327133915399"
$ws.Range("X2").Value = [double]"0.23778664902368299"
$ws.Range("Y2").Value = [double]"0.198778444205298"

# row 3
$ws.Range("B3").Value = [double]"0.19248803450215099"
$ws.Range("C3").Value = [double]"0.13781615317385201"
$ws.Range("D3").Value = [double]"0.207963887580392"
$ws.Range("E3").Value = [double]"0.54650033021872402"
$ws.Range("F3").Value = [double]"0.25109771693949601"
$ws.Range("G3").Value = [double]"3.9622033037627097E-2"
$ws.Range("H3").Value = [double]"0.47580567461714901"
$ws.Range("I3").Value = [double]"0.30249149779601398"
$ws.Range("J3").Value = [double]"0.37457897197872603"
$ws.Range("K3").Value = [double]"0.23801898638354299"
$ws.Range("L3").Value = [double]"0.66237808128884701"
$ws.Range("M3").Value = [double]"9.4207665824658105E-3"
$ws.Range("N3").Value = [double]"0.50897653006355903"
$ws.Range("O3").Value = [double]"0.27209436159616701"
$ws.Range("P3").Value = [double]"0.37745258477374499"
$ws.Range("Q3").Value = [double]"0.125530623063108"
$ws.Range("R3").Value = [double]"0"
$ws.Range("S3").Value = [double]"0.10659888390639501"
$ws.Range("T3").Value = [double]"9.5900803631832399E-2"
$ws.Range("U3").Value = [double]"0.69366396987406798"
$ws.Range("V3").Value = [double]"0.45106932379249898"
$ws.Range("W3").Value = [double]"2.2556717876606801E-2"
$ws.Range("X3").Value = [double]"0.33931273734450401"
$ws.Range("Y3").Value = [double]"1.2238686902251199"

# row 4
$ws.Range("B4").Value = [double]"0.20924603761409499"
$ws.Range("C4").Value = [double]"0.119648048114009"
$ws.Range("D4").Value = [double]"0.20859408180880701"
$ws.Range("E4").Value = [double]"0.54012388947640499"
$ws.Range("F4").Value = [double]"0.272598195191805"
$ws.Range("G4").Value = [double]"4.0676906843209799E-2"
$ws.Range("H4").Value = [double]"0.50880374107016002"
$ws.Range("I4").Value = [double]"0.30464349910744398"
$ws.Range("J4").Value = [double]"0.47683254638237099"
$ws.Range("K4").Value = [double]"0.27739159360244597"
$ws.Range("L4").Value = [double]"1.02013434459466"
$ws.Range("M4").Value = [double]"9.8408047734070102E-3"
$ws.Range("N4").Value = [double]"0.49518670483786398"
$ws.Range("O4").Value = [double]"0"
$ws.Range("P4").Value = [double]"0.13641152181885899"
$ws.Range("Q4").Value = [double]"2.1013061073826299E-2"
$ws.Range("R4").Value = [double]"5.2672724350133697E-2"
$ws.Range("S4").Value = [double]"8.8769204228857507E-2"
$ws.Range("T4").Value = [double]"0"
$ws.Range("U4").Value = [double]"0.170874922345507"
$ws.Range("V4").Value = [double]"0.425319816443819"
$ws.Range("W4").Value = [double]"0.209440406217621"
$ws.Range("X4").Value = [double]"0.198965214903018"
$ws.Range("Y4").Value = [double]"0.66618009473170703"

# row 5
$ws.Range("B5").Value = [double]"0.25790746541027398"
$ws.Range("C5").Value = [double]"0.185513970705079"
$ws.Range("D5").Value = [double]"0.18366486073518901"
$ws.Range("E5").Value = [double]"0.68216320424850496"
$ws.Range("F5").Value = [double]"0.32300188737082602"
$ws.Range("G5").Value = [double]"3.7556935010021497E-2"
$ws.Range("H5").Value = [double]"0.51955878674268896"
$ws.Range("I5").Value = [double]"0.32994338167271398"
$ws.Range("J5").Value = [double]"0.39333507098411002"
$ws.Range("K5").Value = [double]"0.28450716199328702"
$ws.Range("L5").Value = [double]"0.82059427329882695"
$ws.Range("M5").Value = [double]"1.96441672830053E-2"
$ws.Range("N5").Value = [double]"0.48008791353899"
$ws.Range("O5").Value = [double]"0.106882234055029"
$ws.Range("P5").Value = [double]"0.21051312408668099"
$ws.Range("Q5").Value = [double]"0.19967981779500499"
$ws.Range("R5").Value = [double]"9.2380576115702703E-2"
$ws.Range("S5").Value = [double]"0.157497750863339"
$ws.Range("T5").Value = [double]"8.9763650263735104E-2"
$ws.Range("U5").Value = [double]"0.29658804752990697"
$ws.Range("V5").Value = [double]"0.181650577162349"
$ws.Range("W5").Value = [double]"5.7651659725749302E-2"
$ws.Range("X5").Value = [double]"0"
$ws.Range("Y5").Value = [double]"0.58810375111909596"

# row 6
$ws.Range("B6").Value = [double]"0.145830299387498"
$ws.Range("C6").Value = [double]"0.207743282577418"
$ws.Range("D6").Value = [double]"0.16926909537135301"
$ws.Range("E6").Value = [double]"0.52375376479254199"
$ws.Range("F6").Value = [double]"0.18146799070167799"
$ws.Range("G6").Value = [double]"6.2006988333092698E-3"
$ws.Range("H6").Value = [double]"0.283025120825382"
$ws.Range("I6").Value = [double]"0.113561500952091"
$ws.Range("J6").Value = [double]"0.38712908038109201"
$ws.Range("K6").Value = [double]"0.24457253479167801"
$ws.Range("L6").Value = [double]"0.70350145175308099"
$ws.Range("M6").Value = [double]"1.7257648998398201E-2"
$ws.Range("N6").Value = [double]"0.45289417570072898"
$ws.Range("O6").Value = [double]"4.9986101095975399E-2"
$ws.Range("P6").Value = [double]"0.20703343485079401"
$ws.Range("Q6").Value = [double]"0.21657869062681601"
$ws.Range("R6").Value = [double]"0.16910613447006101"
$ws.Range("S6").Value = [double]"0.111469402055943"
$ws.Range("T6").Value = [double]"0.33666199008339098"
$ws.Range("U6").Value = [double]"0.39258375922457001"
$ws.Range("V6").Value = [double]"0.21226604148926201"
$ws.Range("W6").Value = [double]"0.22176034892360599"
$ws.Range("X6").Value = [double]"0.42830935212863103"
$ws.Range("Y6").Value = [double]"0.45121780970434699"

# row 7
$ws.Range("B7").Value = [double]"0"
$ws.Range("C7").Value = [double]"0.46268240540572297"
$ws.Range("D7").Value = [double]"0.135709657210245"
$ws.Range("E7").Value = [double]"0.12532326053963999"
$ws.Range("F7").Value = [double]"0"
$ws.Range("G7").Value = [double]"2.3086579100799701E-2"
$ws.Range("H7").Value = [double]"5.4620532585260201E-2"
$ws.Range("I7").Value = [double]"0"
$ws.Range("J7").Value = [double]"0.35106797701062098"
$ws.Range("K7").Value = [double]"0.148590010946763"
$ws.Range("L7").Value = [double]"0.68664211403178599"
$ws.Range("M7").Value = [double]"0"
$ws.Range("N7").Value = [double]"0.11475017399868"
$ws.Range("O7").Value = [double]"0.15021346084140999"
$ws.Range("P7").Value = [double]"0.201768083140798"
$ws.Range("Q7").Value = [double]"0.62939652749798702"
$ws.Range("R7").Value = [double]"0.20212734478510899"
$ws.Range("S7").Value = [double]"0.204678531020897"
$ws.Range("T7").Value = [double]"0.361117364905991"
$ws.Range("U7").Value = [double]"1.0893274507025501"
$ws.Range("V7").Value = [double]"0.41228016178738203"
$ws.Range("W7").Value = [double]"3.2219558616814199E-2"
$ws.Range("X7").Value = [double]"0.25428325403915297"
$ws.Range("Y7").Value = [double]"1.0656804005840701"

# row 8
$ws.Range("B8").Value = [double]"0.12411971754875201"
$ws.Range("C8").Value = [double]"0.41055408209402899"
$ws.Range("D8").Value = [double]"0.12982640078358401"
$ws.Range("E8").Value = [double]"0.26956901110367898"
$ws.Range("F8").Value = [double]"0.140596425910261"
$ws.Range("G8").Value = [double]"1.5056295796324899E-2"
$ws.Range("H8").Value = [double]"9.2687404891335506E-2"
$ws.Range("I8").Value = [double]"0.20775092846222301"
$ws.Range("J8").Value = [double]"0.32995147683175902"
$ws.Range("K8").Value = [double]"0.244890611114011"
$ws.Range("L8").Value = [double]"0.51524894660522202"
$ws.Range("M8").Value = [double]"1.55826558492567E-2"
$ws.Range("N8").Value = [double]"0.31180034790529099"
$ws.Range("O8").Value = [double]"0.22152510820388199"
$ws.Range("P8").Value = [double]"0.26524741452919898"
$ws.Range("Q8").Value = [double]"0.587457058160195"
$ws.Range("R8").Value = [double]"0.18777494861348301"
$ws.Range("S8").Value = [double]"8.6946342811810795E-2"
$ws.Range("T8").Value = [double]"0.60492334306559203"
$ws.Range("U8").Value = [double]"0.65035203449342105"
$ws.Range("V8").Value = [double]"0"
$ws.Range("W8").Value = [double]"1.40481302709733E-2"
$ws.Range("X8").Value = [double]"2.11462082201762E-2"
$ws.Range("Y8").Value = [double]"0.651668288013174"

# row 9
$ws.Range("B9").Value = [double]"3.4699651349003997E-2"
$ws.Range("C9").Value = [double]"0.29779487560020002"
$ws.Range("D9").Value = [double]"0.16557393354725999"
$ws.Range("E9").Value = [double]"0"
$ws.Range("F9").Value = [double]"7.2284837933941901E-2"
$ws.Range("G9").Value = [double]"1.7763897694751798E-2"
$ws.Range("H9").Value = [double]"0"
$ws.Range("I9").Value = [double]"5.3148920798814003E-2"
$ws.Range("J9").Value = [double]"0.24658733891217599"
$ws.Range("K9").Value = [double]"0.18785586408128099"
$ws.Range("L9").Value = [double]"0.40761829456277099"
$ws.Range("M9").Value = [double]"1.52863536514725E-2"
$ws.Range("N9").Value = [double]"0"
$ws.Range("O9").Value = [double]"0.15434591615197399"
$ws.Range("P9").Value = [double]"0.26246941546617703"
$ws.Range("Q9").Value = [double]"0.44592546307167502"
$ws.Range("R9").Value = [double]"0.13580108921686401"
$ws.Range("S9").Value = [double]"0.15206723599497601"
$ws.Range("T9").Value = [double]"0.41610667057124201"
$ws.Range("U9").Value = [double]"0.60742681218230998"
$ws.Range("V9").Value = [double]"0.68887844067337101"
$ws.Range("W9").Value = [double]"0.157601576973911"
$ws.Range("X9").Value = [double]"0.92654810928600295"
$ws.Range("Y9").Value = [double]"1.75410432061856"

# row 10
$ws.Range("B10").Value = [double]"5.5043167704210597E-2"
$ws.Range("C10").Value = [double]"0.416347544611882"
$ws.Range("D10").Value = [double]"0.106440613291041"
$ws.Range("E10").Value = [double]"3.3678074950925702E-2"
$ws.Range("F10").Value = [double]"6.9558583513097905E-2"
$ws.Range("G10").Value = [double]"1.50298803740455E-2"
$ws.Range("H10").Value = [double]"1.4957401238488899E-2"
$ws.Range("I10").Value = [double]"0.10237413288404"
$ws.Range("J10").Value = [double]"0.22167574166669299"
$ws.Range("K10").Value = [double]"0.12188679621095901"
$ws.Range("L10").Value = [double]"0.65883530583113803"
$ws.Range("M10").Value = [double]"1.94145194039906E-2"
$ws.Range("N10").Value = [double]"8.2043926349042104E-2"
$ws.Range("O10").Value = [double]"0.314543844346098"
$ws.Range("P10").Value = [double]"0.34044144660750097"
$ws.Range("Q10").Value = [double]"0.53127604058754496"
$ws.Range("R10").Value = [double]"0.143250723271758"
$ws.Range("S10").Value = [double]"0.138302568157783"
$ws.Range("T10").Value = [double]"0.50717591490484504"
$ws.Range("U10").Value = [double]"1.72183723727988"
$ws.Range("V10").Value = [double]"0.147780185025227"
$ws.Range("W10").Value = [double]"0"
$ws.Range("X10").Value = [double]"0.28820648691962802"
$ws.Range("Y10").Value = [double]"0"

# row 11
$ws.Range("B11").Value = [double]"0.10678788950474"
$ws.Range("C11").Value = [double]"0.33546254193361702"
$ws.Range("D11").Value = [double]"0.18080738820804501"
$ws.Range("E11").Value = [double]"0.14975547721967999"
$ws.Range("F11").Value = [double]"0.12598262228295901"
$ws.Range("G11").Value = [double]"2.23576211853302E-2"
$ws.Range("H11").Value = [double]"0.10077774067925201"
$ws.Range("I11").Value = [double]"4.2925284247142002E-2"
$ws.Range("J11").Value = [double]"0.225430479902582"
$ws.Range("K11").Value = [double]"0.16927508875042699"
$ws.Range("L11").Value = [double]"0.25098103028144397"
$ws.Range("M11").Value = [double]"1.30127255150123E-5"
$ws.Range("N11").Value = [double]"0.21777288283194399"
$ws.Range("O11").Value = [double]"4.3600471441967499E-2"
$ws.Range("P11").Value = [double]"0"
$ws.Range("Q11").Value = [double]"0.45680866153376498"
$ws.Range("R11").Value = [double]"3.2987304931007003E-2"
$ws.Range("S11").Value = [double]"5.0368704214651699E-2"
$ws.Range("T11").Value = [double]"0.41292524517556201"
$ws.Range("U11").Value = [double]"1.34943396744667"
$ws.Range("V11").Value = [double]"0.26713998539188299"
$ws.Range("W11").Value = [double]"8.8877842938051097E-3"
$ws.Range("X11").Value = [double]"0.359415429938858"
$ws.Range("Y11").Value = [double]"0.89929949923373598"

# row 12
$ws.Range("B12").Value = [double]"0.14062270921536801"
$ws.Range("C12").Value = [double]"0.24938182577788701"
$ws.Range("D12").Value = [double]"0"
$ws.Range("E12").Value = [double]"0.35890923552562998"
$ws.Range("F12").Value = [double]"0.17871573683342601"
$ws.Range("G12").Value = [double]"3.2135613970379899E-3"
$ws.Range("H12").Value = [double]"0.27709721402271198"
$ws.Range("I12").Value = [double]"0.177433495677082"
$ws.Range("J12").Value = [double]"0.175377925132297"
$ws.Range("K12").Value = [double]"0.13299866082696499"
$ws.Range("L12").Value = [double]"1.28734010966686"
$ws.Range("M12").Value = [double]"1.23061015482346E-2"
$ws.Range("N12").Value = [double]"0.32548802722454301"
$ws.Range("O12").Value = [double]"0.14592239136291299"
$ws.Range("P12").Value = [double]"0.234098120144758"
$ws.Range("Q12").Value = [double]"0.31952591687338999"
$ws.Range("R12").Value = [double]"0.10306415962895001"
$ws.Range("S12").Value = [double]"0.11150241055441"
$ws.Range("T12").Value = [double]"0.31422688262512799"
$ws.Range("U12").Value = [double]"0.66009660364017797"
$ws.Range("V12").Value = [double]"0.304211343683297"
$ws.Range("W12").Value = [double]"0.107786189790074"
$ws.Range("X12").Value = [double]"0.27779647101122801"
$ws.Range("Y12").Value = [double]"0.74701054826407098"

# row 13
$ws.Range("B13").Value = [double]"0.14062270921536801"
$ws.Range("C13").Value = [double]"0.24938182577788701"
$ws.Range("D13").Value = [double]"0.217870123618956"
$ws.Range("E13").Value = [double]"0.35890923552562998"
$ws.Range("F13").Value = [double]"0.17871573683342601"
$ws.Range("G13").Value = [double]"0"
$ws.Range("H13").Value = [double]"0.27709721402271198"
$ws.Range("I13").Value = [double]"0.177433495677082"
$ws.Range("J13").Value = [double]"0"
$ws.Range("K13").Value = [double]"2.91139651372028E-2"
$ws.Range("L13").Value = [double]"0"
$ws.Range("M13").Value = [double]"1.49716033300131E-2"
$ws.Range("N13").Value = [double]"0.32548802722454301"
$ws.Range("O13").Value = [double]"0.14592239136291299"
$ws.Range("P13").Value = [double]"0.234098120144758"
$ws.Range("Q13").Value = [double]"0.31952591687338999"
$ws.Range("R13").Value = [double]"0.10306415962895001"
$ws.Range("S13").Value = [double]"0.11150241055441"
$ws.Range("T13").Value = [double]"0.31422688262512799"
$ws.Range("U13").Value = [double]"0.66009660364017797"
$ws.Range("V13").Value = [double]"0.304211343683297"
$ws.Range("W13").Value = [double]"0.107786189790074"
$ws.Range("X13").Value = [double]"0.27779647101122801"
$ws.Range("Y13").Value = [double]"0.74701054826407098"

# row 14
$ws.Range("B14").Value = [double]"0.14062270921536801"
$ws.Range("C14").Value = [double]"0.24938182577788701"
$ws.Range("D14").Value = [double]"0.12734109934544"
$ws.Range("E14").Value = [double]"0.35890923552562998"
$ws.Range("F14").Value = [double]"0.17871573683342601"
$ws.Range("G14").Value = [double]"8.0440101086094798E-3"
$ws.Range("H14").Value = [double]"0.27709721402271198"
$ws.Range("I14").Value = [double]"0.177433495677082"
$ws.Range("J14").Value = [double]"0.19645679974645899"
$ws.Range("K14").Value = [double]"0"
$ws.Range("L14").Value = [double]"0.50075504608704402"
$ws.Range("M14").Value = [double]"5.5493873373917597E-3"
$ws.Range("N14").Value = [double]"0.32548802722454301"
$ws.Range("O14").Value = [double]"0.14592239136291299"
$ws.Range("P14").Value = [double]"0.234098120144758"
$ws.Range("Q14").Value = [double]"0.31952591687338999"
$ws.Range("R14").Value = [double]"0.10306415962895001"
$ws.Range("S14").Value = [double]"0.11150241055441"
$ws.Range("T14").Value = [double]"0.31422688262512799"
$ws.Range("U14").Value = [double]"0.66009660364017797"
$ws.Range("V14").Value = [double]"0.304211343683297"
$ws.Range("W14").Value = [double]"0.107786189790074"
$ws.Range("X14").Value = [double]"0.27779647101122801"
$ws.Range("Y14").Value = [double]"0.74701054826407098"

# row 15
$ws.Range("B15").Value = [double]"0.26165630801588002"
$ws.Range("C15").Value = [double]"0"
$ws.Range("D15").Value = [double]"0.21453380926243501"
$ws.Range("E15").Value = [double]"0.65377136185123097"
$ws.Range("F15").Value = [double]"0.32850609747766701"
$ws.Range("G15").Value = [double]"2.5234502701915702E-2"
$ws.Range("H15").Value = [double]"0.57481168264668303"
$ws.Range("I15").Value = [double]"0.33080510878459701"
$ws.Range("J15").Value = [double]"0.43600030461078298"
$ws.Range("K15").Value = [double]"0.28253992094286801"
$ws.Range("L15").Value = [double]"1.0795756256700999"
$ws.Range("M15").Value = [double]"2.2435640270931499E-2"
$ws.Range("N15").Value = [double]"0.56062549089103897"
$ws.Range("O15").Value = [double]"6.4071858350242802E-2"
$ws.Range("P15").Value = [double]"0.20858483411988199"
$ws.Range("Q15").Value = [double]"0"
$ws.Range("R15").Value = [double]"3.49161887127298E-2"
$ws.Range("S15").Value = [double]"0"
$ws.Range("T15").Value = [double]"0.159400865497001"
$ws.Range("U15").Value = [double]"0"
$ws.Range("V15").Value = [double]"0.30278570641990699"
$ws.Range("W15").Value = [double]"0.21643722165489199"
$ws.Range("X15").Value = [double]"0.14786077722525301"
$ws.Range("Y15").Value = [double]"0.66812602307487701"

# row 16
$ws.Range("B16").Value = [double]"0.16420178758104001"
$ws.Range("C16").Value = [double]"0.17029795241622001"
$ws.Range("D16").Value = [double]"0.209814259383728"
$ws.Range("E16").Value = [double]"0.45179507878800201"
$ws.Range("F16").Value = [double]"0.21332195689522901"
$ws.Range("G16").Value = [double]"2.0349093230160601E-2"
$ws.Range("H16").Value = [double]"0.39608854861225301"
$ws.Range("I16").Value = [double]"0.169366062345608"
$ws.Range("J16").Value = [double]"0.37871042276154998"
$ws.Range("K16").Value = [double]"0.27263441750637402"
$ws.Range("L16").Value = [double]"1.0633808567960199"
$ws.Range("M16").Value = [double]"2.8329976590492099E-2"
$ws.Range("N16").Value = [double]"0.410138348589906"
$ws.Range("O16").Value = [double]"0.13447229653214299"
$ws.Range("P16").Value = [double]"0.26449994673785998"
$ws.Range("Q16").Value = [double]"0.28303289654670799"
$ws.Range("R16").Value = [double]"5.6210341576417001E-2"
$ws.Range("S16").Value = [double]"7.5728807215040994E-2"
$ws.Range("T16").Value = [double]"0.31664175031307901"
$ws.Range("U16").Value = [double]"0.38421527848668602"
$ws.Range("V16").Value = [double]"0.30930077518351801"
$ws.Range("W16").Value = [double]"0.12995471340763001"
$ws.Range("X16").Value = [double]"0.16329863823641699"
$ws.Range("Y16").Value = [double]"0.77398406985333901"

# row 17
$ws.Range("B17").Value = [double]"0.11144766371716699"
$ws.Range("C17").Value = [double]"0.28090096019607202"
$ws.Range("D17").Value = [double]"0.21905931302918899"
$ws.Range("E17").Value = [double]"0.28418213664696801"
$ws.Range("F17").Value = [double]"0.13820413721418701"
$ws.Range("G17").Value = [double]"1.0430347232857899E-2"
$ws.Range("H17").Value = [double]"0.23495726417395801"
$ws.Range("I17").Value = [double]"9.6180049324713496E-2"
$ws.Range("J17").Value = [double]"0.35410137964719302"
$ws.Range("K17").Value = [double]"0.21708260307730501"
$ws.Range("L17").Value = [double]"0.98680588658847601"
$ws.Range("M17").Value = [double]"2.8759990595225299E-2"
$ws.Range("N17").Value = [double]"0.25899113184288403"
$ws.Range("O17").Value = [double]"0.219250187636164"
$ws.Range("P17").Value = [double]"0.27063489074061903"
$ws.Range("Q17").Value = [double]"0.49360766690261598"
$ws.Range("R17").Value = [double]"0.15538742473119199"
$ws.Range("S17").Value = [double]"9.6694049881516003E-2"
$ws.Range("T17").Value = [double]"0.452719562826772"
$ws.Range("U17").Value = [double]"0.58198223990212306"
$ws.Range("V17").Value = [double]"0.28483797689811402"
$ws.Range("W17").Value = [double]"0.16326907797019599"
$ws.Range("X17").Value = [double]"0.24622126588064"
$ws.Range("Y17").Value = [double]"0.77012573606958401"

# Apply scientific-notation number format to M11 (matches new style index 3 / numFmtId 11)
$ws.Range("M11").NumberFormat = "0.00E+00"

# Update sheet dimension implicitly handled by engine; set worksheet as active/selected with L11 selected
$ws.Activate()
$ws.Range("L11").Select()
